$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The production_date column (D) holds text values like "YYYY-MM-DD",
# not real Excel dates. Force Text number format first so assigning a
# "2026-02-.." style string doesn't get auto-parsed into a date serial.
$dateCells = @(
    @{ Cell = "D2"; Value = "2026-02-12" },
    @{ Cell = "D3"; Value = "2026-02-13" },
    @{ Cell = "D4"; Value = "2026-02-14" },
    @{ Cell = "D5"; Value = "2026-02-15" },
    @{ Cell = "D6"; Value = "2026-02-16" }
)

foreach ($entry in $dateCells) {
    $range = $ws.Range($entry.Cell)
    $range.NumberFormat = "@"
    $range.Value = $entry.Value
}
